$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the betting predictions grid for the Spain GP (new race results shuffle
# the P1-P12 / Pole picks for each player). Target values derived from the
# canonical OOXML diff.
$ws.Cells.Item(2,5).Value = "VER"  # E2: LEC -> VER
$ws.Cells.Item(2,6).Value = "NOR"  # F2: VER -> NOR
$ws.Cells.Item(2,7).Value = "LEC"  # G2: HAM -> LEC
$ws.Cells.Item(2,8).Value = "HAM"  # H2: NOR -> HAM
$ws.Cells.Item(2,9).Value = "RUS"  # I2: HAD -> RUS
$ws.Cells.Item(2,10).Value = "ANT"  # J2: LAW -> ANT
$ws.Cells.Item(2,11).Value = "ALO"  # K2: ANT -> ALO
$ws.Cells.Item(2,12).Value = "HAD"  # L2: ALO -> HAD
$ws.Cells.Item(2,13).Value = "GAS"  # M2: RUS -> GAS
$ws.Cells.Item(2,15).Value = "BOR"  # O2: ALB -> BOR
$ws.Cells.Item(3,2).Value = "NOR"  # B3: LEC -> NOR
$ws.Cells.Item(3,4).Value = "VER"  # D3: PIA -> VER
$ws.Cells.Item(3,5).Value = "PIA"  # E3: VER -> PIA
$ws.Cells.Item(3,6).Value = "NOR"  # F3: LEC -> NOR
$ws.Cells.Item(3,7).Value = "ANT"  # G3: NOR -> ANT
$ws.Cells.Item(3,8).Value = "HAM"  # H3: RUS -> HAM
$ws.Cells.Item(3,9).Value = "LEC"  # I3: HAM -> LEC
$ws.Cells.Item(3,10).Value = "RUS"  # J3: ANT -> RUS
$ws.Cells.Item(3,11).Value = "LAW"  # K3: ALO -> LAW
$ws.Cells.Item(3,12).Value = "TSU"  # L3: GAS -> TSU
$ws.Cells.Item(3,13).Value = "ALO"  # M3: OCO -> ALO
$ws.Cells.Item(3,14).Value = "GAS"  # N3: STR -> GAS
$ws.Cells.Item(3,15).Value = "STR"  # O3: SAI -> STR
$ws.Cells.Item(4,2).Value = "PIA"  # B4: LEC -> PIA
$ws.Cells.Item(4,4).Value = "PIA"  # D4: LEC -> PIA
$ws.Cells.Item(4,5).Value = "NOR"  # E4: PIA -> NOR
$ws.Cells.Item(4,6).Value = "VER"  # F4: HAM -> VER
$ws.Cells.Item(4,7).Value = "RUS"  # G4: VER -> RUS
$ws.Cells.Item(4,8).Value = "LEC"  # H4: NOR -> LEC
$ws.Cells.Item(4,9).Value = "ANT"  # I4: RUS -> ANT
$ws.Cells.Item(4,10).Value = "HAM"  # J4: ALB -> HAM
$ws.Cells.Item(4,11).Value = "HAD"  # K4: ANT -> HAD
$ws.Cells.Item(4,12).Value = "ALO"  # L4: TSU -> ALO
$ws.Cells.Item(4,13).Value = "GAS"  # M4: SAI -> GAS
$ws.Cells.Item(4,14).Value = "ALB"  # N4: HAD -> ALB
$ws.Cells.Item(4,15).Value = "BOR"  # O4: ALO -> BOR
$ws.Cells.Item(5,2).Value = "NOR"  # B5: PIA -> NOR
$ws.Cells.Item(5,5).Value = "PIA"  # E5: HAM -> PIA
$ws.Cells.Item(5,6).Value = "VER"  # F5: GAS -> VER
$ws.Cells.Item(5,7).Value = "LEC"  # G5: VER -> LEC
$ws.Cells.Item(5,9).Value = "HAM"  # I5: LEC -> HAM
$ws.Cells.Item(5,10).Value = "ANT"  # J5: HAD -> ANT
$ws.Cells.Item(5,12).Value = "OCO"  # L5: TSU -> OCO
$ws.Cells.Item(5,13).Value = "SAI"  # M5: HAM -> SAI
$ws.Cells.Item(5,14).Value = "STR"  # N5: SAI -> STR
$ws.Cells.Item(5,15).Value = "HAD"  # O5: BEA -> HAD
$ws.Cells.Item(6,2).Value = "PIA"  # B6: LEC -> PIA
$ws.Cells.Item(6,4).Value = "PIA"  # D6: LEC -> PIA
$ws.Cells.Item(6,5).Value = "VER"  # E6: PIA -> VER
$ws.Cells.Item(6,7).Value = "LEC"  # G6: VER -> LEC
$ws.Cells.Item(6,11).Value = "ALO"  # K6: ALB -> ALO
$ws.Cells.Item(6,12).Value = "HAD"  # L6: SAI -> HAD
$ws.Cells.Item(6,13).Value = "GAS"  # M6: ALO -> GAS
$ws.Cells.Item(6,15).Value = "SAI"  # O6: TSU -> SAI
$ws.Cells.Item(7,2).Value = "PIA"  # B7: LEC -> PIA
$ws.Cells.Item(7,4).Value = "PIA"  # D7: LEC -> PIA
$ws.Cells.Item(7,5).Value = "NOR"  # E7: PIA -> NOR
$ws.Cells.Item(7,6).Value = "VER"  # F7: NOR -> VER
$ws.Cells.Item(7,7).Value = "RUS"  # G7: VER -> RUS
$ws.Cells.Item(7,8).Value = "LEC"  # H7: HAM -> LEC
$ws.Cells.Item(7,9).Value = "ANT"  # I7: RUS -> ANT
$ws.Cells.Item(7,10).Value = "HAM"  # J7: ANT -> HAM
$ws.Cells.Item(7,11).Value = "HAD"  # K7: ALB -> HAD
$ws.Cells.Item(7,13).Value = "GAS"  # M7: SAI -> GAS
$ws.Cells.Item(7,14).Value = "SAI"  # N7: TSU -> SAI
$ws.Cells.Item(7,15).Value = "ALB"  # O7: LAW -> ALB
$ws.Cells.Item(8,9).Value = "HAM"  # I8: TSU -> HAM
$ws.Cells.Item(8,10).Value = "TSU"  # J8: HAM -> TSU
$ws.Cells.Item(8,12).Value = "ALO"  # L8: SAI -> ALO
$ws.Cells.Item(8,13).Value = "SAI"  # M8: ALO -> SAI
$ws.Cells.Item(9,2).Value = "PIA"  # B9: LEC -> PIA
$ws.Cells.Item(9,4).Value = "PIA"  # D9: LEC -> PIA
$ws.Cells.Item(9,5).Value = "NOR"  # E9: HAM -> NOR
$ws.Cells.Item(9,6).Value = "VER"  # F9: PIA -> VER
$ws.Cells.Item(9,7).Value = "RUS"  # G9: NOR -> RUS
$ws.Cells.Item(9,8).Value = "LEC"  # H9: VER -> LEC
$ws.Cells.Item(9,9).Value = "HAM"  # I9: ANT -> HAM
$ws.Cells.Item(9,10).Value = "ANT"  # J9: RUS -> ANT
$ws.Cells.Item(9,11).Value = "TSU"  # K9: ALB -> TSU
$ws.Cells.Item(9,12).Value = "ALB"  # L9: TSU -> ALB
$ws.Cells.Item(9,13).Value = "ALO"  # M9: HAD -> ALO
$ws.Cells.Item(9,14).Value = "GAS"  # N9: LAW -> GAS
$ws.Cells.Item(10,2).Value = "PIA"  # B10: LEC -> PIA
$ws.Cells.Item(10,4).Value = "PIA"  # D10: LEC -> PIA
$ws.Cells.Item(10,5).Value = "NOR"  # E10: PIA -> NOR
$ws.Cells.Item(10,7).Value = "RUS"  # G10: HAM -> RUS
$ws.Cells.Item(10,8).Value = "LEC"  # H10: NOR -> LEC
$ws.Cells.Item(10,9).Value = "ANT"  # I10: RUS -> ANT
$ws.Cells.Item(10,10).Value = "HAM"  # J10: ANT -> HAM
$ws.Cells.Item(10,11).Value = "ALO"  # K10: ALB -> ALO
$ws.Cells.Item(10,12).Value = "HAD"  # L10: ALO -> HAD
$ws.Cells.Item(10,13).Value = "GAS"  # M10: TSU -> GAS
$ws.Cells.Item(10,14).Value = "TSU"  # N10: SAI -> TSU
$ws.Cells.Item(10,15).Value = "SAI"  # O10: HAD -> SAI
$ws.Cells.Item(11,2).Value = "PIA"  # B11: LEC -> PIA
$ws.Cells.Item(11,4).Value = "PIA"  # D11: LEC -> PIA
$ws.Cells.Item(11,5).Value = "NOR"  # E11: PIA -> NOR
$ws.Cells.Item(11,6).Value = "VER"  # F11: NOR -> VER
$ws.Cells.Item(11,7).Value = "RUS"  # G11: VER -> RUS
$ws.Cells.Item(11,8).Value = "LEC"  # H11: HAM -> LEC
$ws.Cells.Item(11,9).Value = "ANT"  # I11: RUS -> ANT
$ws.Cells.Item(11,10).Value = "HAM"  # J11: ANT -> HAM
$ws.Cells.Item(11,11).Value = "GAS"  # K11: SAI -> GAS
$ws.Cells.Item(11,12).Value = "HAD"  # L11: ALB -> HAD
$ws.Cells.Item(11,14).Value = "ALO"  # N11: HAD -> ALO
$ws.Cells.Item(11,15).Value = "LAW"  # O11: ALO -> LAW
$ws.Cells.Item(13,5).Value = "NOR"  # E13: VER -> NOR
$ws.Cells.Item(13,6).Value = "VER"  # F13: NOR -> VER
$ws.Cells.Item(13,7).Value = "RUS"  # G13: LEC -> RUS
$ws.Cells.Item(13,8).Value = "LEC"  # H13: HAM -> LEC
$ws.Cells.Item(13,9).Value = "HAM"  # I13: RUS -> HAM
$ws.Cells.Item(13,13).Value = "GAS"  # M13: SAI -> GAS
$ws.Cells.Item(13,14).Value = "SAI"  # N13: GAS -> SAI
$ws.Cells.Item(14,5).Value = "NOR"  # E14: VER -> NOR
$ws.Cells.Item(14,6).Value = "VER"  # F14: NOR -> VER
$ws.Cells.Item(14,7).Value = "RUS"  # G14: LEC -> RUS
$ws.Cells.Item(14,8).Value = "LEC"  # H14: RUS -> LEC
$ws.Cells.Item(14,9).Value = "ANT"  # I14: HAM -> ANT
$ws.Cells.Item(14,10).Value = "HAM"  # J14: TSU -> HAM
$ws.Cells.Item(14,11).Value = "TSU"  # K14: ANT -> TSU
$ws.Cells.Item(14,12).Value = "HAD"  # L14: ALO -> HAD
$ws.Cells.Item(14,13).Value = "GAS"  # M14: SAI -> GAS
$ws.Cells.Item(14,14).Value = "LAW"  # N14: ALB -> LAW
$ws.Cells.Item(14,15).Value = "SAI"  # O14: LAW -> SAI
$ws.Cells.Item(15,6).Value = "RUS"  # F15: VER -> RUS
$ws.Cells.Item(15,7).Value = "VER"  # G15: LEC -> VER
$ws.Cells.Item(15,8).Value = "LEC"  # H15: HAM -> LEC
$ws.Cells.Item(15,9).Value = "ANT"  # I15: RUS -> ANT
$ws.Cells.Item(15,10).Value = "HAM"  # J15: ANT -> HAM
$ws.Cells.Item(15,11).Value = "HAD"  # K15: SAI -> HAD
$ws.Cells.Item(15,13).Value = "SAI"  # M15: ALB -> SAI
$ws.Cells.Item(15,14).Value = "LAW"  # N15: ALO -> LAW
$ws.Cells.Item(15,15).Value = "GAS"  # O15: HAD -> GAS
$ws.Cells.Item(16,2).Value = "NOR"  # B16: LEC -> NOR
$ws.Cells.Item(16,4).Value = "NOR"  # D16: LEC -> NOR
$ws.Cells.Item(16,5).Value = "VER"  # E16: PIA -> VER
$ws.Cells.Item(16,6).Value = "PIA"  # F16: HAM -> PIA
$ws.Cells.Item(16,7).Value = "LEC"  # G16: NOR -> LEC
$ws.Cells.Item(16,8).Value = "RUS"  # H16: VER -> RUS
$ws.Cells.Item(16,9).Value = "HAM"  # I16: RUS -> HAM
$ws.Cells.Item(16,10).Value = "ANT"  # J16: ALB -> ANT
$ws.Cells.Item(16,11).Value = "HAD"  # K16: ANT -> HAD
$ws.Cells.Item(16,12).Value = "LAW"  # L16: SAI -> LAW
$ws.Cells.Item(16,13).Value = "TSU"  # M16: ALO -> TSU
$ws.Cells.Item(16,14).Value = "ALO"  # N16: HAD -> ALO
$ws.Cells.Item(17,9).Value = "HAM"  # I17: ANT -> HAM
$ws.Cells.Item(17,10).Value = "ANT"  # J17: TSU -> ANT
$ws.Cells.Item(17,11).Value = "TSU"  # K17: HAM -> TSU
$ws.Cells.Item(17,14).Value = "ALO"  # N17: DOO -> ALO
$ws.Cells.Item(18,2).Value = "NOR"  # B18: LEC -> NOR
$ws.Cells.Item(18,4).Value = "NOR"  # D18: LEC -> NOR
$ws.Cells.Item(18,6).Value = "VER"  # F18: HAM -> VER
$ws.Cells.Item(18,7).Value = "RUS"  # G18: NOR -> RUS
$ws.Cells.Item(18,8).Value = "LEC"  # H18: VER -> LEC
$ws.Cells.Item(18,9).Value = "HAM"  # I18: RUS -> HAM
$ws.Cells.Item(18,10).Value = "ALB"  # J18: HAD -> ALB
$ws.Cells.Item(18,11).Value = "TSU"  # K18: ALB -> TSU
$ws.Cells.Item(18,12).Value = "ANT"  # L18: LAW -> ANT
$ws.Cells.Item(18,13).Value = "SAI"  # M18: TSU -> SAI
$ws.Cells.Item(18,14).Value = "OCO"  # N18: ALO -> OCO
$ws.Cells.Item(18,15).Value = "HAD"  # O18: SAI -> HAD
$ws.Cells.Item(19,2).Value = "NOR"  # B19: LEC -> NOR
$ws.Cells.Item(19,4).Value = "NOR"  # D19: LEC -> NOR
$ws.Cells.Item(19,7).Value = "RUS"  # G19: NOR -> RUS
$ws.Cells.Item(19,8).Value = "LEC"  # H19: HAM -> LEC
$ws.Cells.Item(19,9).Value = "HAM"  # I19: RUS -> HAM
$ws.Cells.Item(19,10).Value = "ALB"  # J19: HAD -> ALB
$ws.Cells.Item(19,11).Value = "TSU"  # K19: ALB -> TSU
$ws.Cells.Item(19,12).Value = "ANT"  # L19: LAW -> ANT
$ws.Cells.Item(19,13).Value = "SAI"  # M19: TSU -> SAI
$ws.Cells.Item(19,14).Value = "HAD"  # N19: ALO -> HAD
$ws.Cells.Item(19,15).Value = "OCO"  # O19: SAI -> OCO
$ws.Cells.Item(20,2).Value = "PIA"  # B20: VER -> PIA
$ws.Cells.Item(20,4).Value = "PIA"  # D20: VER -> PIA
$ws.Cells.Item(20,5).Value = "VER"  # E20: LEC -> VER
$ws.Cells.Item(20,7).Value = "RUS"  # G20: HAM -> RUS
$ws.Cells.Item(20,8).Value = "ANT"  # H20: PIA -> ANT
$ws.Cells.Item(20,9).Value = "LEC"  # I20: RUS -> LEC
$ws.Cells.Item(20,12).Value = "HAM"  # L20: ANT -> HAM
$ws.Cells.Item(20,14).Value = "ALB"  # N20: TSU -> ALB
$ws.Cells.Item(20,15).Value = "HAD"  # O20: ALB -> HAD
$ws.Cells.Item(21,2).Value = "PIA"  # B21: VER -> PIA
$ws.Cells.Item(21,8).Value = "RUS"  # H21: HAM -> RUS
$ws.Cells.Item(21,9).Value = "HAM"  # I21: ANT -> HAM
$ws.Cells.Item(21,10).Value = "ANT"  # J21: RUS -> ANT
$ws.Cells.Item(21,11).Value = "ALO"  # K21: SAI -> ALO
$ws.Cells.Item(21,12).Value = "ALB"  # L21: STR -> ALB
$ws.Cells.Item(21,13).Value = "GAS"  # M21: OCO -> GAS
$ws.Cells.Item(21,14).Value = "SAI"  # N21: TSU -> SAI
$ws.Cells.Item(21,15).Value = "TSU"  # O21: GAS -> TSU
$ws.Cells.Item(22,9).Value = "HAM"  # I22: SAI -> HAM
$ws.Cells.Item(22,10).Value = "ANT"  # J22: HAM -> ANT
$ws.Cells.Item(22,11).Value = "HAD"  # K22: ANT -> HAD
$ws.Cells.Item(22,12).Value = "ALB"  # L22: TSU -> ALB
$ws.Cells.Item(22,13).Value = "SAI"  # M22: ALB -> SAI
$ws.Cells.Item(22,14).Value = "TSU"  # N22: HAD -> TSU
$ws.Cells.Item(22,15).Value = "ALO"  # O22: LAW -> ALO
$ws.Cells.Item(23,2).Value = "PIA"  # B23: LEC -> PIA
$ws.Cells.Item(23,4).Value = "PIA"  # D23: LEC -> PIA
$ws.Cells.Item(23,5).Value = "NOR"  # E23: PIA -> NOR
$ws.Cells.Item(23,6).Value = "VER"  # F23: HAM -> VER
$ws.Cells.Item(23,7).Value = "RUS"  # G23: NOR -> RUS
$ws.Cells.Item(23,8).Value = "ANT"  # H23: VER -> ANT
$ws.Cells.Item(23,9).Value = "LEC"  # I23: ALB -> LEC
$ws.Cells.Item(23,10).Value = "HAM"  # J23: ANT -> HAM
$ws.Cells.Item(23,11).Value = "ALO"  # K23: RUS -> ALO
$ws.Cells.Item(23,12).Value = "HAD"  # L23: SAI -> HAD
$ws.Cells.Item(23,13).Value = "GAS"  # M23: HAD -> GAS
$ws.Cells.Item(23,14).Value = "SAI"  # N23: ALO -> SAI
$ws.Cells.Item(23,15).Value = "ALB"  # O23: TSU -> ALB
$ws.Cells.Item(24,5).Value = "NOR"  # E24: LEC -> NOR
$ws.Cells.Item(24,7).Value = "LEC"  # G24: NOR -> LEC
$ws.Cells.Item(24,8).Value = "RUS"  # H24: HAM -> RUS
$ws.Cells.Item(24,9).Value = "HAM"  # I24: RUS -> HAM
$ws.Cells.Item(24,10).Value = "ANT"  # J24: SAI -> ANT
$ws.Cells.Item(24,11).Value = "HAD"  # K24: ALO -> HAD
$ws.Cells.Item(24,12).Value = "GAS"  # L24: ALB -> GAS
$ws.Cells.Item(24,13).Value = "SAI"  # M24: ANT -> SAI
$ws.Cells.Item(24,14).Value = "ALB"  # N24: HAD -> ALB
$ws.Cells.Item(24,15).Value = "TSU"  # O24: GAS -> TSU
$ws.Cells.Item(25,5).Value = "NOR"  # E25: LEC -> NOR
$ws.Cells.Item(25,6).Value = "VER"  # F25: NOR -> VER
$ws.Cells.Item(25,7).Value = "RUS"  # G25: VER -> RUS
$ws.Cells.Item(25,8).Value = "LEC"  # H25: HAM -> LEC
$ws.Cells.Item(25,9).Value = "HAM"  # I25: RUS -> HAM
$ws.Cells.Item(25,13).Value = "TSU"  # M25: HAD -> TSU
$ws.Cells.Item(25,14).Value = "HAD"  # N25: GAS -> HAD
$ws.Cells.Item(25,15).Value = "LAW"  # O25: ALO -> LAW
$ws.Cells.Item(26,2).Value = "NOR"  # B26: LEC -> NOR
$ws.Cells.Item(26,4).Value = "NOR"  # D26: LEC -> NOR
$ws.Cells.Item(26,5).Value = "VER"  # E26: NOR -> VER
$ws.Cells.Item(26,7).Value = "SAI"  # G26: VER -> SAI
$ws.Cells.Item(26,10).Value = "ANT"  # J26: SAI -> ANT
$ws.Cells.Item(26,13).Value = "LAW"  # M26: ANT -> LAW
$ws.Cells.Item(26,14).Value = "OCO"  # N26: STR -> OCO
$ws.Cells.Item(26,15).Value = "HUL"  # O26: TSU -> HUL

# Move the active selection to F29 (matches the saved workbook view state).
$ws.Range("F29").Select()
